$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110; existing rows 110..212 shift down to 111..213.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new observation.
$ws.Range("A110").Value = 3
$ws.Range("B110").Value = "Femacal de La Calera"
$ws.Range("C110").Value = "Coquimbo"
$ws.Range("D110").Value = 44484
$ws.Range("E110").Value = 5
$ws.Range("F110").Value = 100112043
$ws.Range("G110").Value = "Pepino ensalada"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 73
$ws.Range("K110").Value = 13000
$ws.Range("L110").Value = 14000
$ws.Range("M110").Value = 13479
$ws.Range("N110").Value = "$/caja 70 unidades"
$ws.Range("O110").Value = "Región de Arica y Parinacota"
$ws.Range("P110").Value = 193
$ws.Range("Q110").Value = 70
$ws.Range("R110").Value = "Hortaliza"
